$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.76"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.33%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.373"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.17%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08368"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.939"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9709"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.62%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.523"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.41%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1124"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.02%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1896"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.39%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.32%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04608"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.35%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1062"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.36%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001306"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.64%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005888"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.93%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.367"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.09%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.416"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.52%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.16%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.566"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-16.85%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1388"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.27%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2574"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.31%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04155"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.14%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001237"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.06%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004404"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.18%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.52%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002980"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.35%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02716"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-2.32%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05623"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.82%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007817"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.27%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1413"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.61%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007316"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.03%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002040"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.34%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008697"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.51%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3506"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006890"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003486"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-6.29%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003532"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.07%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
